$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "242.97"
$ws.Cells.Item(2, 7).Value = "18"

# Row 3
$ws.Cells.Item(3, 4).Value = "23.02"
$ws.Cells.Item(3, 7).Value = "18"

# Row 4
$ws.Cells.Item(4, 2).Value = "LEO"
$ws.Cells.Item(4, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(4, 4).Value = "3.617"
$ws.Cells.Item(4, 5).Value = "3LEOLEO"
$ws.Cells.Item(4, 7).Value = "18"

# Row 5
$ws.Cells.Item(5, 2).Value = "HuobiToken"
$ws.Cells.Item(5, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(5, 4).Value = "5.397"
$ws.Cells.Item(5, 5).Value = "4HuobiTokenHT"
$ws.Cells.Item(5, 7).Value = "18"

# Row 6
$ws.Cells.Item(6, 2).Value = "Cronos"
$ws.Cells.Item(6, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.05930"
$ws.Cells.Item(6, 5).Value = "5CronosCRO"
$ws.Cells.Item(6, 7).Value = "18"

# Row 7
$ws.Cells.Item(7, 2).Value = "GateToken"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(7, 4).Value = "3.393"
$ws.Cells.Item(7, 5).Value = "6GateTokenGT"
$ws.Cells.Item(7, 7).Value = "18"

# Row 8
$ws.Cells.Item(8, 2).Value = "KuCoinToken"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Cells.Item(8, 4).Value = "6.457"
$ws.Cells.Item(8, 5).Value = "7KuCoinTokenKCS"
$ws.Cells.Item(8, 7).Value = "18"

# Row 9
$ws.Cells.Item(9, 2).Value = "MXToken"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(9, 4).Value = "0.8074"
$ws.Cells.Item(9, 5).Value = "8MXTokenMX"
$ws.Cells.Item(9, 7).Value = "18"

# Row 10
$ws.Cells.Item(10, 2).Value = "FTXToken"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(10, 4).Value = "0.9076"
$ws.Cells.Item(10, 5).Value = "9FTXTokenFTT"
$ws.Cells.Item(10, 7).Value = "18"

# Row 11
$ws.Cells.Item(11, 2).Value = "One"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Cells.Item(11, 4).Value = "0.01108"
$ws.Cells.Item(11, 5).Value = "10OneONE"
$ws.Cells.Item(11, 7).Value = "18"

# Row 12
$ws.Cells.Item(12, 2).Value = "WazirX"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(12, 4).Value = "0.1416"
$ws.Cells.Item(12, 5).Value = "11WazirXWRX"
$ws.Cells.Item(12, 7).Value = "18"

# Row 13
$ws.Cells.Item(13, 2).Value = "MandalaExchangeToken"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(13, 4).Value = "0.07444"
$ws.Cells.Item(13, 5).Value = "12MandalaExchangeTokenMDX"
$ws.Cells.Item(13, 7).Value = "18"

# Row 14
$ws.Cells.Item(14, 2).Value = "LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(14, 4).Value = "0.03317"
$ws.Cells.Item(14, 5).Value = "13LiechtensteinCryptoassetsExchangeLCX"
$ws.Cells.Item(14, 7).Value = "18"

# Row 15
$ws.Cells.Item(15, 2).Value = "BitrueCoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(15, 4).Value = "0.03044"
$ws.Cells.Item(15, 5).Value = "14BitrueCoinBTR"
$ws.Cells.Item(15, 7).Value = "18"

# Row 16
$ws.Cells.Item(16, 2).Value = "BitMartToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(16, 4).Value = "0.09325"
$ws.Cells.Item(16, 5).Value = "15BitMartTokenBMX"
$ws.Cells.Item(16, 7).Value = "18"

# Row 17
$ws.Cells.Item(17, 2).Value = "MCDex"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Cells.Item(17, 4).Value = "3.951"
$ws.Cells.Item(17, 5).Value = "16MCDexMCB"
$ws.Cells.Item(17, 7).Value = "18"

# Row 18
$ws.Cells.Item(18, 2).Value = "BitForexToken"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(18, 4).Value = "0.001583"
$ws.Cells.Item(18, 5).Value = "17BitForexTokenBF"
$ws.Cells.Item(18, 7).Value = "18"

# Row 19
$ws.Cells.Item(19, 2).Value = "CoinExToken"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Cells.Item(19, 4).Value = "0.04802"
$ws.Cells.Item(19, 5).Value = "18CoinExTokenCET"
$ws.Cells.Item(19, 7).Value = "18"

# Row 20
$ws.Cells.Item(20, 2).Value = "TigerCash"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(20, 4).Value = "0.006088"
$ws.Cells.Item(20, 5).Value = "19TigerCashTCH"
$ws.Cells.Item(20, 7).Value = "18"

# Row 21
$ws.Cells.Item(21, 2).Value = "UpBots"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Cells.Item(21, 4).Value = "0.007493"
$ws.Cells.Item(21, 5).Value = "20UpBotsUBXTBestin24h"
$ws.Cells.Item(21, 7).Value = "18"

# Row 22
$ws.Cells.Item(22, 2).Value = "HotbitToken"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Cells.Item(22, 4).Value = "0.004415"
$ws.Cells.Item(22, 5).Value = "21HotbitTokenHTB"
$ws.Cells.Item(22, 7).Value = "18"

# Row 23
$ws.Cells.Item(23, 2).Value = "BitKan"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Cells.Item(23, 4).Value = "0.0009864"
$ws.Cells.Item(23, 5).Value = "22BitKanKAN"
$ws.Cells.Item(23, 7).Value = "18"

# Row 24
$ws.Cells.Item(24, 2).Value = "NitroEx"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Cells.Item(24, 4).Value = "0.00007805"
$ws.Cells.Item(24, 5).Value = "23NitroExNTX"
$ws.Cells.Item(24, 7).Value = "18"

# Row 25
$ws.Cells.Item(25, 7).Value = "18"

# Row 26
$ws.Cells.Item(26, 7).Value = "18"

# Row 27
$ws.Cells.Item(27, 7).Value = "18"

# Row 28
$ws.Cells.Item(28, 7).Value = "18"

# Row 29
$ws.Cells.Item(29, 7).Value = "18"

# Row 30
$ws.Cells.Item(30, 7).Value = "18"

# Row 31
$ws.Cells.Item(31, 7).Value = "18"

# Row 32
$ws.Cells.Item(32, 7).Value = "18"

# Row 33
$ws.Cells.Item(33, 7).Value = "18"

# Row 34
$ws.Cells.Item(34, 7).Value = "18"

# Row 35
$ws.Cells.Item(35, 7).Value = "18"

# Row 36
$ws.Cells.Item(36, 7).Value = "18"

# Row 37
$ws.Cells.Item(37, 7).Value = "18"

# Row 38
$ws.Cells.Item(38, 7).Value = "18"

# Row 39
$ws.Cells.Item(39, 7).Value = "18"

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.03870"
$ws.Cells.Item(40, 7).Value = "18"

# Row 41
$ws.Cells.Item(41, 4).Value = "0.006195"
$ws.Cells.Item(41, 7).Value = "18"

# Row 42
$ws.Cells.Item(42, 4).Value = "0.1064"
$ws.Cells.Item(42, 7).Value = "18"

# Row 43
$ws.Cells.Item(43, 7).Value = "18"

# Row 44
$ws.Cells.Item(44, 4).Value = "0.007245"
$ws.Cells.Item(44, 7).Value = "18"

# Row 45
$ws.Cells.Item(45, 4).Value = "0.00005186"
$ws.Cells.Item(45, 7).Value = "18"

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.00000000750"
$ws.Cells.Item(46, 7).Value = "18"

# Row 47
$ws.Cells.Item(47, 4).Value = "0.0005806"
$ws.Cells.Item(47, 7).Value = "18"

# Row 48
$ws.Cells.Item(48, 4).Value = "0.9006"
$ws.Cells.Item(48, 7).Value = "18"

# Row 49
$ws.Cells.Item(49, 7).Value = "18"

# Row 50
$ws.Cells.Item(50, 4).Value = "0.00002101"
$ws.Cells.Item(50, 7).Value = "18"

# Row 51
$ws.Cells.Item(51, 4).Value = "0.0002001"
$ws.Cells.Item(51, 7).Value = "18"
